# Añadida función de enviar el resultado al correo y se han escondido las variables sensibles
#
# Appends the two new Google News rows that were scraped after the last run
# to the bottom of the "news" sheet (rows 6 and 7), growing the used range
# from A1:C5 to A1:C7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{
        published = "05-11-2024"
        title     = "El Gobierno invierte un millón de euros para la mejora del abastecimiento de agua en el norte de Gran Canaria - InfonorteDigital"
        link      = "https://news.google.com/rss/articles/CBMi4wFBVV95cUxOYkNvMmRBd0tBamZONVNxbzdtYy1ZR3l6SWlQTHFaTUI5ajJaZ190bGVaYkRKcWFLV0xQS0owRGdVcFlMSU1qTC1CVFZnWS1UdWlidEZqeDRnTG1yZGNXOXF2NHVBVGU2UzJhNUNTMFhIWTQ2WVFNMWQ5RmVDTWNfRkJ2WVppTWVGaFNZNHFMRnRqNVdWZnZjMWREZlk1RGdRVW9CMW8ySHlVbkwwUERjanhvMF8zdVl3UHloeExUZXJDWi1kZllxR0M2eGJCczNHeFpwRmxkMU9HNXRjQVBPY0R2TQ?oc=5"
    },
    @{
        published = "06-11-2024"
        title     = "La Feria Internacional Canagua y Energía impulsa un espacio para los encuentros profesionales entre empresas y visitantes - Canarias Noticias"
        link      = "https://news.google.com/rss/articles/CBMiugFBVV95cUxOMXMzRFBiNnRZUmNrY1JrVmN3S24zSjFPVjBEcVIwTEg4czZORnhIZ0ZDVTJnZVl5QlVrT2lMLTQzVUlxeDVQVmtSZVBXX0R1TUZJUjFPcjU0c1BiQV9GOWNPcGFueUEwTGdmakZqZHNEbkJ2YlNoOG9mbHZlUklRX2QyQjBuUUdzeFA1ZmgtaHBDRzNjZVVhUUp5Umtja1B4bGpEMXhYU0NSVDB0TWpHY3RYQnRKaC1nN1E?oc=5"
    }
)

$startRow = 6
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $dateCell = $ws.Cells.Item($r, 1)

    # The "published" column holds dd-mm-yyyy text (e.g. "05-11-2024"). Excel's
    # automatic type detection would otherwise silently reinterpret this as a
    # date serial (since the day portion is <= 12), unlike the existing rows
    # above whose day-of-month is always > 12. Temporarily force the cell to
    # Text so the literal string is preserved, then drop the formatting again
    # so the new rows end up styled the same as the rest of the data (no
    # explicit number format), exactly like rows 2-5.
    $dateCell.NumberFormat = "@"
    $dateCell.Value2 = $row.published
    $dateCell.ClearFormats()

    $ws.Cells.Item($r, 2).Value2 = $row.title
    $ws.Cells.Item($r, 3).Value2 = $row.link
}

Write-Host "Added $($newRows.Count) rows starting at row $startRow"
